# Regenerated "centroides" (KMeans cluster centroids) table - rows 2-11 of
# Sheet1 (columns A:J) are overwritten in place with the updated centroid
# data that comes from adding the modelo_kmeans.predict(...) example and
# re-running the notebook (cluster index <-> row assignment shifted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2  (cluster 0)
$ws.Cells.Item(2, 1).Value  = 0
$ws.Cells.Item(2, 2).Value  = 0
$ws.Cells.Item(2, 3).Value  = 1
$ws.Cells.Item(2, 4).Value  = 0
$ws.Cells.Item(2, 5).Value  = 0
$ws.Cells.Item(2, 6).Value  = 1
$ws.Cells.Item(2, 7).Value  = 0
$ws.Cells.Item(2, 8).Value  = 0
$ws.Cells.Item(2, 9).Value  = 0.5845
$ws.Cells.Item(2, 10).Value = 46.1755

# Row 3  (cluster 1)
$ws.Cells.Item(3, 1).Value  = 1
$ws.Cells.Item(3, 2).Value  = 1
$ws.Cells.Item(3, 3).Value  = 0
$ws.Cells.Item(3, 4).Value  = 1
$ws.Cells.Item(3, 5).Value  = 0
$ws.Cells.Item(3, 6).Value  = 0
$ws.Cells.Item(3, 7).Value  = 0
$ws.Cells.Item(3, 8).Value  = 0
$ws.Cells.Item(3, 9).Value  = 0.5425
$ws.Cells.Item(3, 10).Value = 42.8575

# Row 4  (cluster 2)
$ws.Cells.Item(4, 1).Value  = 2
$ws.Cells.Item(4, 2).Value  = 0
$ws.Cells.Item(4, 3).Value  = 1
$ws.Cells.Item(4, 4).Value  = 1
$ws.Cells.Item(4, 5).Value  = 0
$ws.Cells.Item(4, 6).Value  = 0
$ws.Cells.Item(4, 7).Value  = 0
$ws.Cells.Item(4, 8).Value  = 0
$ws.Cells.Item(4, 9).Value  = 0.5604
$ws.Cells.Item(4, 10).Value = 44.2716

# Row 5  (cluster 3)
$ws.Cells.Item(5, 1).Value  = 3
$ws.Cells.Item(5, 2).Value  = 0
$ws.Cells.Item(5, 3).Value  = 1
$ws.Cells.Item(5, 4).Value  = 0
$ws.Cells.Item(5, 5).Value  = 0
$ws.Cells.Item(5, 6).Value  = 0
$ws.Cells.Item(5, 7).Value  = 1
$ws.Cells.Item(5, 8).Value  = 0
$ws.Cells.Item(5, 9).Value  = 0.5513
$ws.Cells.Item(5, 10).Value = 43.5527

# Row 6  (cluster 4)
$ws.Cells.Item(6, 1).Value  = 4
$ws.Cells.Item(6, 2).Value  = 0
$ws.Cells.Item(6, 3).Value  = 1
$ws.Cells.Item(6, 4).Value  = 0
$ws.Cells.Item(6, 5).Value  = 1
$ws.Cells.Item(6, 6).Value  = 0
$ws.Cells.Item(6, 7).Value  = 0
$ws.Cells.Item(6, 8).Value  = 0
$ws.Cells.Item(6, 9).Value  = 0.5713
$ws.Cells.Item(6, 10).Value = 45.1327

# Row 7  (cluster 5)
$ws.Cells.Item(7, 1).Value  = 5
$ws.Cells.Item(7, 2).Value  = 1
$ws.Cells.Item(7, 3).Value  = 0
$ws.Cells.Item(7, 4).Value  = 0
$ws.Cells.Item(7, 5).Value  = 1
$ws.Cells.Item(7, 6).Value  = 0
$ws.Cells.Item(7, 7).Value  = 0
$ws.Cells.Item(7, 8).Value  = 0
$ws.Cells.Item(7, 9).Value  = 0.5537
$ws.Cells.Item(7, 10).Value = 43.7423

# Row 8  (cluster 6)
$ws.Cells.Item(8, 1).Value  = 6
$ws.Cells.Item(8, 2).Value  = 0
$ws.Cells.Item(8, 3).Value  = 1
$ws.Cells.Item(8, 4).Value  = 0
$ws.Cells.Item(8, 5).Value  = 0
$ws.Cells.Item(8, 6).Value  = 0
$ws.Cells.Item(8, 7).Value  = 0
$ws.Cells.Item(8, 8).Value  = 1
$ws.Cells.Item(8, 9).Value  = 0.5339
$ws.Cells.Item(8, 10).Value = 42.1781

# Row 9  (cluster 7)
$ws.Cells.Item(9, 1).Value  = 7
$ws.Cells.Item(9, 2).Value  = 1
$ws.Cells.Item(9, 3).Value  = 0
$ws.Cells.Item(9, 4).Value  = 0
$ws.Cells.Item(9, 5).Value  = 0
$ws.Cells.Item(9, 6).Value  = 0
$ws.Cells.Item(9, 7).Value  = 0
$ws.Cells.Item(9, 8).Value  = 1
$ws.Cells.Item(9, 9).Value  = 0.5661
$ws.Cells.Item(9, 10).Value = 44.72190000000001

# Row 10 (cluster 8)
$ws.Cells.Item(10, 1).Value  = 8
$ws.Cells.Item(10, 2).Value  = 1
$ws.Cells.Item(10, 3).Value  = 0
$ws.Cells.Item(10, 4).Value  = 0
$ws.Cells.Item(10, 5).Value  = 0
$ws.Cells.Item(10, 6).Value  = 1
$ws.Cells.Item(10, 7).Value  = 0
$ws.Cells.Item(10, 8).Value  = 0
$ws.Cells.Item(10, 9).Value  = 0.5737
$ws.Cells.Item(10, 10).Value = 45.3223

# Row 11 (cluster 9)
$ws.Cells.Item(11, 1).Value  = 9
$ws.Cells.Item(11, 2).Value  = 1
$ws.Cells.Item(11, 3).Value  = 0
$ws.Cells.Item(11, 4).Value  = 0
$ws.Cells.Item(11, 5).Value  = 0
$ws.Cells.Item(11, 6).Value  = 0
$ws.Cells.Item(11, 7).Value  = 1
$ws.Cells.Item(11, 8).Value  = 0
$ws.Cells.Item(11, 9).Value  = 0.5456
$ws.Cells.Item(11, 10).Value = 43.1024
